# Auto-generated: apply market-data refresh values to computed columns H:N
# across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 315.70456
$ws.Range("I6").Value = 77.125
$ws.Range("J6").Value = 368.72223
$ws.Range("K6").Value = 231.375
$ws.Range("L6").Value = 1106.16669
$ws.Range("M6").Value = -119.375
$ws.Range("N6").Value = -1330.16669
$ws.Range("H13").Value = 10998
$ws.Range("J13").Value = 14996
$ws.Range("L13").Value = 14996
$ws.Range("N13").Value = -15334
$ws.Range("H17").Value = 174478.1
$ws.Range("J17").Value = 174478.1
$ws.Range("L17").Value = 523434.3
$ws.Range("N17").Value = -523770.3
$ws.Range("H40").Value = 3686
$ws.Range("I40").Value = 1999.2
$ws.Range("J40").Value = 4452.727
$ws.Range("K40").Value = 1999.2
$ws.Range("L40").Value = 4452.727
$ws.Range("M40").Value = -1824.2
$ws.Range("N40").Value = -4802.727
$ws.Range("H41").Value = 956.0909
$ws.Range("I41").Value = 466.16666
$ws.Range("J41").Value = 1544
$ws.Range("K41").Value = 466.16666
$ws.Range("L41").Value = 1544
$ws.Range("M41").Value = -26.16665999999998
$ws.Range("N41").Value = -2424
$ws.Range("H80").Value = 1391.6207
$ws.Range("I80").Value = 846.4
$ws.Range("J80").Value = 1678.579
$ws.Range("K80").Value = 2539.2
$ws.Range("L80").Value = 5035.737
$ws.Range("M80").Value = -1541.2
$ws.Range("N80").Value = -7031.737
$ws.Range("H83").Value = 1391.6207
$ws.Range("I83").Value = 846.4
$ws.Range("J83").Value = 1678.579
$ws.Range("K83").Value = 7617.599999999999
$ws.Range("L83").Value = 15107.211
$ws.Range("M83").Value = -2625.599999999999
$ws.Range("N83").Value = -25091.211
$ws.Range("H124").Value = 84999.336
$ws.Range("J124").Value = 84999.336
$ws.Range("L124").Value = 84999.336
$ws.Range("N124").Value = -94819.336
$ws.Range("H138").Value = 441457.7
$ws.Range("I138").Value = 67863.13
$ws.Range("J138").Value = 1001849.5
$ws.Range("K138").Value = 203589.39
$ws.Range("L138").Value = 3005548.5
$ws.Range("M138").Value = -198449.39
$ws.Range("N138").Value = -3015828.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4012.0688
$ws.Range("I61").Value = 3938
$ws.Range("J61").Value = 4475
$ws.Range("K61").Value = 3938
$ws.Range("L61").Value = 4475
$ws.Range("M61").Value = -3726
$ws.Range("N61").Value = -4899
$ws.Range("H74").Value = 2057.6592
$ws.Range("I74").Value = 1722.5405
$ws.Range("J74").Value = 3829
$ws.Range("K74").Value = 1722.5405
$ws.Range("L74").Value = 3829
$ws.Range("M74").Value = -848.5405000000001
$ws.Range("N74").Value = -5577
$ws.Range("H77").Value = 2057.6592
$ws.Range("I77").Value = 1722.5405
$ws.Range("J77").Value = 3829
$ws.Range("K77").Value = 8612.702499999999
$ws.Range("L77").Value = 19145
$ws.Range("M77").Value = -4244.702499999999
$ws.Range("N77").Value = -27881
$ws.Range("H97").Value = 1431.6
$ws.Range("I97").Value = 1431.6
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1431.6
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -935.5999999999999
$ws.Range("N97").ClearContents()
$ws.Range("H132").Value = 1571.3877
$ws.Range("I132").Value = 1277.0227
$ws.Range("J132").Value = 4161.8
$ws.Range("K132").Value = 3831.0681
$ws.Range("L132").Value = 12485.4
$ws.Range("M132").Value = -1301.0681
$ws.Range("N132").Value = -17545.4
$ws.Range("H134").Value = 101177.22
$ws.Range("J134").Value = 101324.375
$ws.Range("L134").Value = 101324.375
$ws.Range("N134").Value = -111464.375
$ws.Range("H136").Value = 4012.0688
$ws.Range("I136").Value = 3938
$ws.Range("J136").Value = 4475
$ws.Range("K136").Value = 11814
$ws.Range("L136").Value = 13425
$ws.Range("M136").Value = -9264
$ws.Range("N136").Value = -18525

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1641.5676
$ws.Range("I20").Value = 1420.375
$ws.Range("J20").Value = 2049.923
$ws.Range("K20").Value = 1420.375
$ws.Range("L20").Value = 2049.923
$ws.Range("M20").Value = -1173.375
$ws.Range("N20").Value = -2543.923
$ws.Range("H86").Value = 10533.214
$ws.Range("I86").Value = 18540.715
$ws.Range("J86").Value = 2525.7144
$ws.Range("K86").Value = 18540.715
$ws.Range("L86").Value = 2525.7144
$ws.Range("M86").Value = -17417.715
$ws.Range("N86").Value = -4771.7144
$ws.Range("H89").Value = 10533.214
$ws.Range("I89").Value = 18540.715
$ws.Range("J89").Value = 2525.7144
$ws.Range("K89").Value = 92703.575
$ws.Range("L89").Value = 12628.572
$ws.Range("M89").Value = -87087.575
$ws.Range("N89").Value = -23860.572
$ws.Range("H94").Value = 1109.1818
$ws.Range("I94").Value = 510.2
$ws.Range("J94").Value = 1608.3334
$ws.Range("K94").Value = 510.2
$ws.Range("L94").Value = 1608.3334
$ws.Range("M94").Value = -59.19999999999999
$ws.Range("N94").Value = -2510.3334
$ws.Range("H134").Value = 1864.6842
$ws.Range("I134").Value = 1140.5667
$ws.Range("J134").Value = 4580.125
$ws.Range("K134").Value = 3421.7001
$ws.Range("L134").Value = 13740.375
$ws.Range("M134").Value = -886.7001
$ws.Range("N134").Value = -18810.375

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 2313.7856
$ws.Range("I5").Value = 733.1667
$ws.Range("J5").Value = 3499.25
$ws.Range("K5").Value = 733.1667
$ws.Range("L5").Value = 3499.25
$ws.Range("M5").Value = -621.1667
$ws.Range("N5").Value = -3723.25
$ws.Range("H31").Value = 2212.244
$ws.Range("I31").Value = 1297.125
$ws.Range("J31").Value = 5466
$ws.Range("K31").Value = 1297.125
$ws.Range("L31").Value = 5466
$ws.Range("M31").Value = -1002.125
$ws.Range("N31").Value = -6056
$ws.Range("H34").Value = 2212.244
$ws.Range("I34").Value = 1297.125
$ws.Range("J34").Value = 5466
$ws.Range("K34").Value = 1297.125
$ws.Range("L34").Value = 5466
$ws.Range("M34").Value = -1095.125
$ws.Range("N34").Value = -5870

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 500022.5
$ws.Range("I9").Value = 666683.3
$ws.Range("J9").Value = 40
$ws.Range("K9").Value = 2000049.9
$ws.Range("L9").Value = 120
$ws.Range("M9").Value = -1999825.9
$ws.Range("N9").Value = -568
$ws.Range("H104").Value = 5316.5
$ws.Range("J104").Value = 5316.5
$ws.Range("L104").Value = 15949.5
$ws.Range("N104").Value = -21191.5
$ws.Range("H131").Value = 1839.5483
$ws.Range("I131").Value = 946.6667
$ws.Range("K131").Value = 2840.0001
$ws.Range("M131").Value = 2199.9999
$ws.Range("H133").Value = 5699
$ws.Range("I133").Value = 1000
$ws.Range("K133").Value = 3000
$ws.Range("M133").Value = 2060

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 347.57144
$ws.Range("I2").Value = 343.77777
$ws.Range("J2").Value = 354.4
$ws.Range("K2").Value = 343.77777
$ws.Range("L2").Value = 354.4
$ws.Range("M2").Value = -230.77777
$ws.Range("N2").Value = -580.4
$ws.Range("H70").Value = 6005
$ws.Range("J70").Value = 6200
$ws.Range("L70").Value = 6200
$ws.Range("N70").Value = -6740
$ws.Range("H73").Value = 6005
$ws.Range("J73").Value = 6200
$ws.Range("L73").Value = 6200
$ws.Range("N73").Value = -8072
$ws.Range("H97").Value = 3382
$ws.Range("I97").Value = 1466.6666
$ws.Range("K97").Value = 1466.6666
$ws.Range("M97").Value = -970.6666
$ws.Range("H126").Value = 2802.4285
$ws.Range("I126").Value = 1874
$ws.Range("K126").Value = 5622
$ws.Range("M126").Value = -3152

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 2069.4614
$ws.Range("I9").Value = 1309.3636
$ws.Range("J9").Value = 6250
$ws.Range("K9").Value = 1309.3636
$ws.Range("L9").Value = 6250
$ws.Range("M9").Value = -1085.3636
$ws.Range("N9").Value = -6698
$ws.Range("H46").Value = 3309.375
$ws.Range("J46").Value = 3681.25
$ws.Range("L46").Value = 3681.25
$ws.Range("N46").Value = -4057.25
$ws.Range("H132").Value = 3113.3171
$ws.Range("I132").Value = 3088.5151
$ws.Range("J132").Value = 3215.625
$ws.Range("K132").Value = 9265.5453
$ws.Range("L132").Value = 9646.875
$ws.Range("M132").Value = -6735.5453
$ws.Range("N132").Value = -14706.875
$ws.Range("H136").Value = 3914.6453
$ws.Range("I136").Value = 3634.56
$ws.Range("J136").Value = 5081.6665
$ws.Range("K136").Value = 10903.68
$ws.Range("L136").Value = 15244.9995
$ws.Range("M136").Value = -8353.68
$ws.Range("N136").Value = -20344.9995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 50172
$ws.Range("J95").Value = 50172
$ws.Range("L95").Value = 50172
$ws.Range("N95").Value = -55664
$ws.Range("H96").Value = 3313.2856
$ws.Range("I96").Value = 1864.3334
$ws.Range("K96").Value = 1864.3334
$ws.Range("M96").Value = -491.3334
$ws.Range("H126").Value = 4112.5
$ws.Range("I126").Value = 2200
$ws.Range("K126").Value = 6600
$ws.Range("M126").Value = -4130
$ws.Range("H132").Value = 2112.8728
$ws.Range("I132").Value = 2072.9565
$ws.Range("J132").Value = 2316.889
$ws.Range("K132").Value = 6218.869499999999
$ws.Range("L132").Value = 6950.667
$ws.Range("M132").Value = -3688.869499999999
$ws.Range("N132").Value = -12010.667
$ws.Range("H136").Value = 10935.381
$ws.Range("I136").Value = 24222.666
$ws.Range("K136").Value = 72667.99800000001
$ws.Range("M136").Value = -70117.99800000001
